$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 286.5
$ws.Range("I28").Value = 299.66666
$ws.Range("J28").Value = 231.2
$ws.Range("K28").Value = 299.66666
$ws.Range("L28").Value = 231.2
$ws.Range("M28").Value = 185.33334
$ws.Range("N28").Value = -1201.2
$ws.Range("H58").Value = 53966.74
$ws.Range("I58").Value = 242.8
$ws.Range("J58").Value = 113660
$ws.Range("K58").Value = 728.4000000000001
$ws.Range("L58").Value = 340980
$ws.Range("M58").Value = -578.4000000000001
$ws.Range("N58").Value = -341280
$ws.Range("H74").Value = 3806.5
$ws.Range("I74").Value = 3143.2856
$ws.Range("J74").Value = 4116
$ws.Range("K74").Value = 3143.2856
$ws.Range("L74").Value = 4116
$ws.Range("M74").Value = -2207.2856
$ws.Range("N74").Value = -5988
$ws.Range("H77").Value = 3806.5
$ws.Range("I77").Value = 3143.2856
$ws.Range("J77").Value = 4116
$ws.Range("K77").Value = 15716.428
$ws.Range("L77").Value = 20580
$ws.Range("M77").Value = -11036.428
$ws.Range("N77").Value = -29940
$ws.Range("H123").Value = 37498.5
$ws.Range("J123").Value = 37498.5
$ws.Range("L123").Value = 37498.5
$ws.Range("N123").Value = -47298.5
$ws.Range("H132").Value = 4654693
$ws.Range("I132").Value = 6064185.5
$ws.Range("K132").Value = 18192556.5
$ws.Range("M132").Value = -18190026.5
$ws.Range("H137").Value = 3639.244
$ws.Range("I137").Value = 3909.3635
$ws.Range("K137").Value = 11728.0905
$ws.Range("M137").Value = -9178.0905
$ws.Range("H138").Value = 3535.86
$ws.Range("I138").Value = 1658.4445
$ws.Range("J138").Value = 5739.7827
$ws.Range("K138").Value = 4975.333500000001
$ws.Range("L138").Value = 17219.3481
$ws.Range("M138").Value = 164.6664999999994
$ws.Range("N138").Value = -27499.3481

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3567.32
$ws.Range("I32").Value = 2845.7559
$ws.Range("J32").Value = 7999.7856
$ws.Range("K32").Value = 2845.7559
$ws.Range("L32").Value = 7999.7856
$ws.Range("M32").Value = -2558.7559
$ws.Range("N32").Value = -8573.785599999999
$ws.Range("H45").Value = 1332.3182
$ws.Range("I45").Value = 1053.2941
$ws.Range("J45").Value = 2281
$ws.Range("K45").Value = 1053.2941
$ws.Range("L45").Value = 2281
$ws.Range("M45").Value = -676.2941000000001
$ws.Range("N45").Value = -3035
$ws.Range("H61").Value = 2820.923
$ws.Range("I61").Value = 1109
$ws.Range("K61").Value = 1109
$ws.Range("M61").Value = -897
$ws.Range("H97").Value = 826.5
$ws.Range("I97").Value = 528.6667
$ws.Range("K97").Value = 528.6667
$ws.Range("M97").Value = -32.66669999999999
$ws.Range("H122").Value = 3179.4666
$ws.Range("I122").Value = 2190.2222
$ws.Range("J122").Value = 4663.3335
$ws.Range("K122").Value = 6570.6666
$ws.Range("L122").Value = 13990.0005
$ws.Range("M122").Value = -4120.6666
$ws.Range("N122").Value = -18890.0005
$ws.Range("H132").Value = 3236.577
$ws.Range("I132").Value = 2242.0588
$ws.Range("J132").Value = 5115.1113
$ws.Range("K132").Value = 6726.176399999999
$ws.Range("L132").Value = 15345.3339
$ws.Range("M132").Value = -4196.176399999999
$ws.Range("N132").Value = -20405.3339
$ws.Range("H136").Value = 2820.923
$ws.Range("I136").Value = 1109
$ws.Range("K136").Value = 3327
$ws.Range("M136").Value = -777

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 18253.484
$ws.Range("I86").Value = 1235.75
$ws.Range("J86").Value = 76600
$ws.Range("K86").Value = 1235.75
$ws.Range("L86").Value = 76600
$ws.Range("M86").Value = -112.75
$ws.Range("N86").Value = -78846
$ws.Range("H89").Value = 18253.484
$ws.Range("I89").Value = 1235.75
$ws.Range("J89").Value = 76600
$ws.Range("K89").Value = 6178.75
$ws.Range("L89").Value = 383000
$ws.Range("M89").Value = -562.75
$ws.Range("N89").Value = -394232
$ws.Range("H94").Value = 512.5
$ws.Range("J94").Value = 387.5
$ws.Range("L94").Value = 387.5
$ws.Range("N94").Value = -1289.5
$ws.Range("H99").Value = 1991.2413
$ws.Range("I99").Value = 1507.9048
$ws.Range("J99").Value = 3260
$ws.Range("K99").Value = 1507.9048
$ws.Range("L99").Value = 3260
$ws.Range("M99").Value = -9.904800000000023
$ws.Range("N99").Value = -6256
$ws.Range("H107").Value = 1904.7819
$ws.Range("I107").Value = 1414.0244
$ws.Range("J107").Value = 3342
$ws.Range("K107").Value = 1414.0244
$ws.Range("L107").Value = 3342
$ws.Range("M107").Value = 505.9756
$ws.Range("N107").Value = -7182
$ws.Range("H134").Value = 6959.533
$ws.Range("I134").Value = 19966.666
$ws.Range("K134").Value = 59899.99800000001
$ws.Range("M134").Value = -57364.99800000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1383.4
$ws.Range("I22").Value = 595.1
$ws.Range("J22").Value = 2960
$ws.Range("K22").Value = 595.1
$ws.Range("L22").Value = 2960
$ws.Range("M22").Value = -245.1
$ws.Range("N22").Value = -3660
$ws.Range("H31").Value = 2932.1738
$ws.Range("I31").Value = 1748.0286
$ws.Range("J31").Value = 6699.909
$ws.Range("K31").Value = 1748.0286
$ws.Range("L31").Value = 6699.909
$ws.Range("M31").Value = -1453.0286
$ws.Range("N31").Value = -7289.909
$ws.Range("H34").Value = 2932.1738
$ws.Range("I34").Value = 1748.0286
$ws.Range("J34").Value = 6699.909
$ws.Range("K34").Value = 1748.0286
$ws.Range("L34").Value = 6699.909
$ws.Range("M34").Value = -1546.0286
$ws.Range("N34").Value = -7103.909
$ws.Range("H58").Value = 9618390
$ws.Range("I58").Value = 1934.7028
$ws.Range("J58").Value = 33338980
$ws.Range("K58").Value = 1934.7028
$ws.Range("L58").Value = 33338980
$ws.Range("M58").Value = -1731.7028
$ws.Range("N58").Value = -33339386
$ws.Range("H134").Value = 1677.92
$ws.Range("I134").Value = 1132.2565
$ws.Range("J134").Value = 3612.5454
$ws.Range("K134").Value = 3396.7695
$ws.Range("L134").Value = 10837.6362
$ws.Range("M134").Value = -861.7694999999999
$ws.Range("N134").Value = -15907.6362
$ws.Range("H136").Value = 9618390
$ws.Range("I136").Value = 1934.7028
$ws.Range("J136").Value = 33338980
$ws.Range("K136").Value = 5804.1084
$ws.Range("L136").Value = 100016940
$ws.Range("M136").Value = -3254.1084
$ws.Range("N136").Value = -100022040

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 15961.25
$ws.Range("J63").Value = 15961.25
$ws.Range("L63").Value = 47883.75
$ws.Range("N63").Value = -49381.75
$ws.Range("H66").Value = 15961.25
$ws.Range("J66").Value = 15961.25
$ws.Range("L66").Value = 143651.25
$ws.Range("N66").Value = -151139.25
$ws.Range("H70").Value = 2644.5715
$ws.Range("I70").Value = 1878
$ws.Range("J70").Value = 3666.6667
$ws.Range("K70").Value = 5634
$ws.Range("L70").Value = 11000.0001
$ws.Range("M70").Value = -5319
$ws.Range("N70").Value = -11630.0001
$ws.Range("H73").Value = 2644.5715
$ws.Range("I73").Value = 1878
$ws.Range("J73").Value = 3666.6667
$ws.Range("K73").Value = 5634
$ws.Range("L73").Value = 11000.0001
$ws.Range("M73").Value = -4542
$ws.Range("N73").Value = -13184.0001
$ws.Range("H87").Value = 8021.3184
$ws.Range("J87").Value = 14587.5
$ws.Range("L87").Value = 43762.5
$ws.Range("N87").Value = -46258.5
$ws.Range("H90").Value = 8021.3184
$ws.Range("J90").Value = 14587.5
$ws.Range("L90").Value = 131287.5
$ws.Range("N90").Value = -143767.5
$ws.Range("H120").Value = 18895.889
$ws.Range("J120").Value = 19004.715
$ws.Range("L120").Value = 57014.145
$ws.Range("N120").Value = -66690.145
$ws.Range("H121").Value = 28731
$ws.Range("I121").Value = 325
$ws.Range("J121").Value = 35043.445
$ws.Range("K121").Value = 975
$ws.Range("L121").Value = 105130.335
$ws.Range("M121").Value = 335
$ws.Range("N121").Value = -107750.335
$ws.Range("H124").Value = 21200
$ws.Range("J124").Value = 100000
$ws.Range("L124").Value = 300000
$ws.Range("N124").Value = -309820

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 29347.027
$ws.Range("I102").Value = 1889.05
$ws.Range("J102").Value = 61650.53
$ws.Range("K102").Value = 1889.05
$ws.Range("L102").Value = 61650.53
$ws.Range("M102").Value = -267.05
$ws.Range("N102").Value = -64894.53
$ws.Range("H113").Value = 2862.7856
$ws.Range("I113").Value = 2688
$ws.Range("J113").Value = 3299.75
$ws.Range("K113").Value = 2688
$ws.Range("L113").Value = 3299.75
$ws.Range("M113").Value = -518
$ws.Range("N113").Value = -7639.75
$ws.Range("H122").Value = 3946.3809
$ws.Range("I122").Value = 2097.5
$ws.Range("J122").Value = 5627.1816
$ws.Range("K122").Value = 6292.5
$ws.Range("L122").Value = 16881.5448
$ws.Range("M122").Value = -3842.5
$ws.Range("N122").Value = -21781.5448
$ws.Range("H126").Value = 3655.8948
$ws.Range("I126").Value = 1940.2222
$ws.Range("J126").Value = 5200
$ws.Range("K126").Value = 5820.6666
$ws.Range("L126").Value = 15600
$ws.Range("M126").Value = -3350.6666
$ws.Range("N126").Value = -20540
$ws.Range("H132").Value = 3054.0784
$ws.Range("I132").Value = 2803.842
$ws.Range("K132").Value = 8411.526
$ws.Range("M132").Value = -5881.526

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1447.6364
$ws.Range("I46").Value = 533.4
$ws.Range("J46").Value = 1716.5294
$ws.Range("K46").Value = 533.4
$ws.Range("L46").Value = 1716.5294
$ws.Range("M46").Value = -345.4
$ws.Range("N46").Value = -2092.5294
$ws.Range("H93").Value = 10000
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 10000
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 10000
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = -12496

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6828.2085
$ws.Range("I132").Value = 1518.762
$ws.Range("J132").Value = 16964.424
$ws.Range("K132").Value = 4556.286
$ws.Range("L132").Value = 50893.272
$ws.Range("M132").Value = -2026.286
$ws.Range("N132").Value = -55953.272
$ws.Range("H136").Value = 978.4194
$ws.Range("I136").Value = 679.05
$ws.Range("K136").Value = 2037.15
$ws.Range("M136").Value = 512.8500000000001
